$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the formatting of the 2020 column (Q) into the new 2021 column (R)
#    for the whole data block (rows 4-44), matching how the author likely
#    extended the table by duplicating the previous year's column.
$ws.Range("Q4:Q44").Copy() | Out-Null
$ws.Range("R4:R44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2. Header: 2021
$ws.Range("R4").Value = 2021

# 3. Data rows (numeric values) per region/indicator block
$ws.Range("R7").Value = 1931.83
$ws.Range("R8").Value = 739818.5

$ws.Range("R11").Value = 1552.9
$ws.Range("R12").Value = 25048.6

$ws.Range("R15").Value = 125.7
$ws.Range("R16").Value = 82213.899999999994

$ws.Range("R19").Value = 99.6
$ws.Range("R20").Value = 80059.600000000006

$ws.Range("R23").Value = 0.9
$ws.Range("R24").Value = 17172.7

$ws.Range("R27").Value = 15.9
$ws.Range("R28").Value = 56666.5

$ws.Range("R31").Value = 58.5
$ws.Range("R32").Value = 30765.1

$ws.Range("R35").Value = 78.3
$ws.Range("R36").Value = 110267.1

$ws.Range("R40").Value = 297797.59999999998

$ws.Range("R44").Value = 39827.4

# 4. "Share" formula rows (percentage of damage = loss / GDP * 100)
$ws.Range("R6").Formula = "=R7/R8*100"
$ws.Range("R10").Formula = "=R11/R12*100"
$ws.Range("R14").Formula = "=R15/R16*100"
$ws.Range("R18").Formula = "=R19/R20*100"
$ws.Range("R22").Formula = "=R23/R24*100"
$ws.Range("R26").Formula = "=R27/R28*100"
$ws.Range("R30").Formula = "=R31/R32*100"
$ws.Range("R34").Formula = "=R35/R36*100"

# 5. Rows with no recorded disasters this year -> "-" placeholder (same text
#    already used elsewhere in the sheet for missing data)
$ws.Range("R38").Value = "-"
$ws.Range("R39").Value = "-"
$ws.Range("R42").Value = "-"
$ws.Range("R43").Value = "-"

# 6. Rows 5, 9, 13, 17, 21, 25, 29, 33, 37, 41 are section headers with no
#    data cell of their own (they only picked up the formatting above).

# 7. Update the view: scroll position and active selection the author left
#    the workbook in after finishing the edit.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("T9").Select() | Out-Null
